# Remove the stray "2014 / Bahrain" row that was produced by the old JSON
# ingestion pipeline. The new CSV-based ingestion no longer emits this row,
# so delete it from the worksheet (shifting all following rows up by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = $null
$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($i = 1; $i -le $rowCount; $i++) {
    $yearVal = $ws.Cells.Item($i, 1).Value()
    $countryVal = $ws.Cells.Item($i, 2).Value()
    if ($yearVal -eq "2014" -and $countryVal -eq "Bahrain") {
        $targetRow = $i
        break
    }
}

if ($targetRow -ne $null) {
    $ws.Rows.Item($targetRow).Delete()
}
